$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.282.08'
$ws.Range("E2").Value = '  +0.10%  '
$ws.Range("D3").Value = '1.867.50'
$ws.Range("E3").Value = '  +0.33%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.70'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.66%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4698'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.35%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2854'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.30%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06568'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.31%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.35'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.59%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07822'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.58%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '96.70'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.96%  '
$ws.Range("D13").Value = '1.835.58'
$ws.Range("E13").Value = '  -1.57%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6956'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.22%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.083'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.16%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '268.43'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.53%  '
$ws.Range("D17").Value = '30.306.04'
$ws.Range("E17").Value = '  +0.19%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.77'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.64%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007686'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.71%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.000'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.10%  '
$ws.Range("D21").Value = '2.116.86'
$ws.Range("E21").Value = '  -0.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.09%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.246'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.54%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.156'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.22%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.561'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.07%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '166.47'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.72%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.86'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.00%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.937'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.44%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.363'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.39%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09903'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.33%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.355'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.14%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.457'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.89%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.048'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.83%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04731'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.61%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.130'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.25%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7034'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.37%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.717'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.30%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01874'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.09%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.770'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.55%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.330'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.24%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '72.78'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.35%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.948'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.82%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4166'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.21%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.000'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.02%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8357'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.58%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '103.05'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.23%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '971.56'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.78%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.112'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.52%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.169'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.00%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.52'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.09%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05684'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.35%  '
